# Remove the "gem conflict" clause from the "Optimized real-time Volt app..." bullet
# under the Gramercy Consultants experience entry.
#
# Before: "...data/model structure; solved gem conflicts by manually specifying a
#          lower version for a dependency; finished features such as reset password..."
# After:  "...data/model structure; finished features such as reset password..."

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "; solved gem conflicts by manually specifying a lower version for a dependency",
    $true,   # MatchCase
    $false,  # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap (wdFindContinue)
    $false,  # Format
    "",      # Replacement text (delete the clause)
    2        # Replace (wdReplaceAll)
)
